# Daily attendance processing - 2026-01-02 11:00:45
# Reorders the "Recorded By" (column G) contributor lists on the
# "Session Analysis Results" worksheet to reflect the latest processing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range('G2').Value = 'system, backup@backdoor.com, System'
$ws.Range('G4').Value = 'backup@backdoor.com, System'
$ws.Range('G5').Value = 'backup@backdoor.com, System'
$ws.Range('G7').Value = 'admin@admin.com, System'
$ws.Range('G8').Value = 'backup@backdoor.com, System'
$ws.Range('G11').Value = 'System, dnasr281@gmail.com'
$ws.Range('G17').Value = 'System, dnasr281@gmail.com'
$ws.Range('G28').Value = 'system, backup@backdoor.com, System'
$ws.Range('G30').Value = 'backup@backdoor.com, System'
$ws.Range('G31').Value = 'backup@backdoor.com, System'
$ws.Range('G33').Value = 'admin@admin.com, System'
$ws.Range('G34').Value = 'backup@backdoor.com, System'
$ws.Range('G37').Value = 'System, dnasr281@gmail.com'
$ws.Range('G43').Value = 'System, dnasr281@gmail.com'
$ws.Range('G54').Value = 'system, backup@backdoor.com, System'
$ws.Range('G56').Value = 'backup@backdoor.com, System'
$ws.Range('G57').Value = 'backup@backdoor.com, System'
$ws.Range('G59').Value = 'admin@admin.com, System'
$ws.Range('G60').Value = 'backup@backdoor.com, System'
$ws.Range('G63').Value = 'System, dnasr281@gmail.com'
$ws.Range('G69').Value = 'System, dnasr281@gmail.com'
$ws.Range('G80').Value = 'backup@backdoor.com, System'
$ws.Range('G81').Value = 'backup@backdoor.com, System'
$ws.Range('G82').Value = 'backup@backdoor.com, System'
$ws.Range('G87').Value = 'admin@admin.com, dnasr281@gmail.com'
$ws.Range('G93').Value = 'System, dnasr281@gmail.com'
$ws.Range('G94').Value = 'System, dnasr281@gmail.com'
$ws.Range('G96').Value = 'System, dnasr281@gmail.com'
$ws.Range('G106').Value = 'backup@backdoor.com, System'
$ws.Range('G107').Value = 'backup@backdoor.com, System'
$ws.Range('G108').Value = 'backup@backdoor.com, System'
$ws.Range('G113').Value = 'admin@admin.com, dnasr281@gmail.com'
$ws.Range('G119').Value = 'System, dnasr281@gmail.com'
$ws.Range('G120').Value = 'System, dnasr281@gmail.com'
$ws.Range('G122').Value = 'System, dnasr281@gmail.com'
$ws.Range('G132').Value = 'backup@backdoor.com, System'
$ws.Range('G133').Value = 'backup@backdoor.com, System'
$ws.Range('G134').Value = 'backup@backdoor.com, System'
$ws.Range('G139').Value = 'admin@admin.com, dnasr281@gmail.com'
$ws.Range('G145').Value = 'System, dnasr281@gmail.com'
$ws.Range('G146').Value = 'System, dnasr281@gmail.com'
$ws.Range('G148').Value = 'System, dnasr281@gmail.com'
